$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 190.58333
$ws.Range("I33").Value = 133.7
$ws.Range("K33").Value = 133.7
$ws.Range("M33").Value = 95.30000000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 844.37036
$ws.Range("I53").Value = 980.13043
$ws.Range("J53").Value = 63.75
$ws.Range("K53").Value = 980.13043
$ws.Range("L53").Value = 63.75
$ws.Range("M53").Value = -343.13043
$ws.Range("N53").Value = -1337.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4127.3687
$ws.Range("I64").Value = 4235.4546
$ws.Range("J64").Value = 3978.75
$ws.Range("K64").Value = 4235.4546
$ws.Range("L64").Value = 3978.75
$ws.Range("M64").Value = -3987.4546
$ws.Range("N64").Value = -4474.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4127.3687
$ws.Range("I67").Value = 4235.4546
$ws.Range("J67").Value = 3978.75
$ws.Range("K67").Value = 4235.4546
$ws.Range("L67").Value = 3978.75
$ws.Range("M67").Value = -3377.4546
$ws.Range("N67").Value = -5694.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1920.138
$ws.Range("I113").Value = 1680.1538
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 1680.1538
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = 1573.8462
$ws.Range("N113").Value = -10508

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2606.875
$ws.Range("I116").Value = 2723.9167
$ws.Range("K116").Value = 2723.9167
$ws.Range("M116").Value = 718.0832999999998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 9013682
$ws.Range("I132").Value = 12826137
$ws.Range("J132").Value = 2426.4546
$ws.Range("K132").Value = 38478411
$ws.Range("L132").Value = 7279.3638
$ws.Range("M132").Value = -38475881
$ws.Range("N132").Value = -12339.3638

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1342.84
$ws.Range("J138").Value = 1633.7273
$ws.Range("L138").Value = 4901.1819
$ws.Range("N138").Value = -15181.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3970.0645
$ws.Range("I32").Value = 3620.309
$ws.Range("K32").Value = 3620.309
$ws.Range("M32").Value = -3333.309

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1152
$ws.Range("I61").Value = 1017.2
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 1017.2
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -805.2
$ws.Range("N61").Value = -2924

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2356.5
$ws.Range("I132").Value = 2337.6365
$ws.Range("J132").Value = 2372.4614
$ws.Range("K132").Value = 7012.9095
$ws.Range("L132").Value = 7117.3842
$ws.Range("M132").Value = -4482.9095
$ws.Range("N132").Value = -12177.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1152
$ws.Range("I136").Value = 1017.2
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 3051.6
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -501.6000000000004
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3369.4092
$ws.Range("I86").Value = 3993.353
$ws.Range("K86").Value = 3993.353
$ws.Range("M86").Value = -2870.353

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3369.4092
$ws.Range("I89").Value = 3993.353
$ws.Range("K89").Value = 19966.765
$ws.Range("M89").Value = -14350.765

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 27778716
$ws.Range("I99").Value = 31250868
$ws.Range("J99").Value = 1499.5
$ws.Range("K99").Value = 31250868
$ws.Range("L99").Value = 1499.5
$ws.Range("M99").Value = -31249370
$ws.Range("N99").Value = -4495.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1494.1578
$ws.Range("I107").Value = 1236.6666
$ws.Range("J107").Value = 2459.75
$ws.Range("K107").Value = 1236.6666
$ws.Range("L107").Value = 2459.75
$ws.Range("M107").Value = 683.3334
$ws.Range("N107").Value = -6299.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4336.9487
$ws.Range("I134").Value = 1084.4517
$ws.Range("J134").Value = 16940.375
$ws.Range("K134").Value = 3253.3551
$ws.Range("L134").Value = 50821.125
$ws.Range("M134").Value = -718.3551000000002
$ws.Range("N134").Value = -55891.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1371.6451
$ws.Range("I31").Value = 917.4375
$ws.Range("J31").Value = 1856.1333
$ws.Range("K31").Value = 917.4375
$ws.Range("L31").Value = 1856.1333
$ws.Range("M31").Value = -622.4375
$ws.Range("N31").Value = -2446.1333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1371.6451
$ws.Range("I34").Value = 917.4375
$ws.Range("J34").Value = 1856.1333
$ws.Range("K34").Value = 917.4375
$ws.Range("L34").Value = 1856.1333
$ws.Range("M34").Value = -715.4375
$ws.Range("N34").Value = -2260.1333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1135.6
$ws.Range("J58").Value = 1489.3334
$ws.Range("L58").Value = 1489.3334
$ws.Range("N58").Value = -1895.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7433276
$ws.Range("I86").Value = 33335636
$ws.Range("K86").Value = 33335636
$ws.Range("M86").Value = -33334513

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 7433276
$ws.Range("I89").Value = 33335636
$ws.Range("K89").Value = 166678180
$ws.Range("M89").Value = -166672564

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1904.5
$ws.Range("I99").Value = 1999.75
$ws.Range("K99").Value = 1999.75
$ws.Range("M99").Value = -501.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 916.2
$ws.Range("I122").Value = 910.3333
$ws.Range("J122").Value = 925
$ws.Range("K122").Value = 2730.9999
$ws.Range("L122").Value = 2775
$ws.Range("M122").Value = -280.9998999999998
$ws.Range("N122").Value = -7675

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1904.5
$ws.Range("I126").Value = 1999.75
$ws.Range("K126").Value = 5999.25
$ws.Range("M126").Value = -3529.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1462.6177
$ws.Range("I132").Value = 836.5769
$ws.Range("J132").Value = 3497.25
$ws.Range("K132").Value = 2509.7307
$ws.Range("L132").Value = 10491.75
$ws.Range("M132").Value = 20.26929999999993
$ws.Range("N132").Value = -15551.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 867.1429000000001
$ws.Range("I134").Value = 779.8333
$ws.Range("J134").Value = 932.625
$ws.Range("K134").Value = 2339.4999
$ws.Range("L134").Value = 2797.875
$ws.Range("M134").Value = 195.5001000000002
$ws.Range("N134").Value = -7867.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1135.6
$ws.Range("J136").Value = 1489.3334
$ws.Range("L136").Value = 4468.0002
$ws.Range("N136").Value = -9568.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 669.25714
$ws.Range("I113").Value = 600
$ws.Range("J113").Value = 671.2941
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 2013.8823
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -6353.882299999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2058.8965
$ws.Range("I132").Value = 1290.5264
$ws.Range("J132").Value = 3518.8
$ws.Range("K132").Value = 3871.5792
$ws.Range("L132").Value = 10556.4
$ws.Range("M132").Value = -1341.5792
$ws.Range("N132").Value = -15616.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 16000
$ws.Range("J64").Value = 16000
$ws.Range("L64").Value = 16000
$ws.Range("N64").Value = -16496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 16000
$ws.Range("J67").Value = 16000
$ws.Range("L67").Value = 16000
$ws.Range("N67").Value = -17716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 358.45456
$ws.Range("I113").Value = 180.5
$ws.Range("J113").Value = 572
$ws.Range("K113").Value = 541.5
$ws.Range("L113").Value = 1716
$ws.Range("M113").Value = 1628.5
$ws.Range("N113").Value = -6056

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 58825050
$ws.Range("I126").Value = 83334940
$ws.Range("J126").Value = 1301
$ws.Range("K126").Value = 250004820
$ws.Range("L126").Value = 3903
$ws.Range("M126").Value = -250002350
$ws.Range("N126").Value = -8843

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1603.16
$ws.Range("I132").Value = 1261.238
$ws.Range("K132").Value = 3783.714
$ws.Range("M132").Value = -1253.714
